# 自动更新Excel文件
# For every data row (row 2 onward): column D = total days, E = remaining
# days, F = start date (yyyymmdd, numeric). Each run represents one day
# passing: the "remaining" count is decremented by one. When a row's
# remaining count has reached 1 (i.e. it is about to expire), the cycle is
# renewed: remaining is reset back to the total (D), and the start date is
# pushed forward by 7 days. Rows whose date in F cannot be parsed as a
# valid yyyymmdd date are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $D = $dCell.Value2
    $E = $eCell.Value2
    $F = $fCell.Value2

    if ($D -eq $null -or $E -eq $null -or $F -eq $null) {
        continue
    }

    $Fstr = [string]([int64]$F)

    $parsedDate = $null
    if ($Fstr.Length -eq 8) {
        try {
            $y = [int]$Fstr.Substring(0, 4)
            $m = [int]$Fstr.Substring(4, 2)
            $day = [int]$Fstr.Substring(6, 2)
            $parsedDate = Get-Date -Year $y -Month $m -Day $day
        } catch {
            $parsedDate = $null
        }
    }

    if ($parsedDate -eq $null) {
        # Unparseable date (e.g. malformed entry) - skip this row entirely.
        continue
    }

    $remaining = [int]$E
    $total = [int]$D

    if ($remaining -eq 1) {
        # Cycle expired - renew it: reset remaining days and roll the
        # start date forward by one week.
        $newRemaining = $total
        $newDate = $parsedDate.AddDays(7)
        $newFstr = $newDate.ToString("yyyyMMdd")

        $eCell.Value = $newRemaining
        $fCell.Value = [int64]$newFstr
    } else {
        $newRemaining = $remaining - 1
        $eCell.Value = $newRemaining
    }
}
